$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New rows of data (row, date-serial, unidad_medida text, cantidad, fc)
$data = @(
    @(44101, "LT", 63561, 0),
    @(44102, "%", 359053, 0),
    @(44103, "UN", 12, 0),
    @(44104, "UN", 10, 0),
    @(44105, "UN", 30, 0),
    @(44106, "M3", 1786, 0),
    @(44107, "M", 280, 0),
    @(44108, "UN", 2, 0),
    @(44109, "m2", 956, 0),
    @(44110, "m3", 919, 0),
    @(44111, "m3", 2434, 0),
    @(44112, "m3", 2757, 0),
    @(44113, "m3", 63, 0),
    @(44114, "M", 46, 0),
    @(44115, "m3", 1344, 0),
    @(44116, "kg", 26797, 0),
    @(44117, "m2", 3432, 0),
    @(44118, "m3", 799, 0),
    @(44119, "m3", 96, 0),
    @(44120, "M2", 3432, 0),
    @(44121, "GL", 3119, 0),
    @(44122, "m3", 347, 0),
    @(44123, "KG", 1021, 0),
    @(44124, "KG", 1021, 0),
    @(44125, "HM", 1, 0),
    @(44126, "LT", 1097, 0),
    @(44127, "m3", 848, 0),
    @(44128, "un", 52, 0),
    @(44129, "un", 12, 0),
    @(44130, "un", 12, 0),
    @(44131, "un", 1184, 0),
    @(44132, "un", 46, 0),
    @(44133, "un", 442, 0),
    @(44134, "un", 89, 0),
    @(44135, "un", 104, 0),
    @(44136, "m3", 1736, 0),
    @(44137, "UN", 158, 0),
    @(44138, "VIAJE", 2, 0)
)

$startRow = 28
$row = $startRow
foreach ($item in $data) {
    $ws.Cells.Item($row, 1).Value = $item[0]
    $ws.Cells.Item($row, 2).Value = $item[1]
    $ws.Cells.Item($row, 3).Value = $item[2]
    $ws.Cells.Item($row, 4).Value = $item[3]
    $row++
}
$endRow = $row - 1

# Column A carries the alternating banded-fill style used by the existing
# rows (odd source row 26 = fill "3", even source row 27 = fill "4"); copy
# that two-row formatting pattern down across all of the newly added rows.
[void]$ws.Range("A26:A27").Copy()
[void]$ws.Range(("A{0}:A{1}" -f $startRow, $endRow)).PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Match the saved selection / scroll position from the edited workbook.
$excel.ActiveWindow.ScrollRow = 48
$excel.ActiveWindow.ScrollColumn = 1
[void]$ws.Range("C66").Select()
